$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ADAM): MBM_Selected 3 -> 5
$ws.Range("C2").Value = 5

# Row 3 (PAM): MBM_Selected 0 -> 3, Working TRUE -> FALSE
$ws.Range("C3").Value = 3
$ws.Range("F3").Value = $false

# Row 4 (MIKE): MBM_Selected 2 -> 1, Working TRUE -> FALSE
$ws.Range("C4").Value = 1
$ws.Range("F4").Value = $false

# Row 5 (CHRIS): MBM_Selected 1 -> 4
$ws.Range("C5").Value = 4

# Row 6 (ANTHONY): MBM_Worked 0 -> 1, MBM_Selected 5 -> 2
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 2

# Row 7 (DAWNETTA): MBM_Worked 0 -> 1, MBM_Selected 4 -> 0
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0
